$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 48
$ws.Cells.Item($row, 1).Value = "I magnifici 2.0"
$ws.Cells.Item($row, 2).Value = "Stefano  Tita | Clitoriders"
$ws.Cells.Item($row, 3).Value = "Alessio Bragagna | SHARK ATTACK"
$ws.Cells.Item($row, 4).Value = "Pietro  Gasparini | Mai una gioia"
$ws.Cells.Item($row, 5).Value = "Geremia  Carollo | FC SAVIGNANO"
$ws.Cells.Item($row, 6).Value = "Mattia Tezzele | U.SGUARNA"
